# Update the weekly schedule for Aline S. M.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value  = "-"
$ws.Range("E6").Value  = "['MEC-1B-T. M. Metalicos', 'MEC-1B-T. M. Metalicos']"
$ws.Range("F10").Value = "[-, -, 'MEC-2A-Metalografia', -]"
$ws.Range("E11").Value = "-"
$ws.Range("C12").Value = "['MEC-1A-T. M. Metalicos', 'MEC-1A-T. M. Metalicos']"
$ws.Range("E12").Value = "-"
$ws.Range("E14").Value = "-"
$ws.Range("F14").Value = "[-, 'MEC-2A-Metalografia', -, -]"
$ws.Range("C15").Value = "-"
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "[-, 'MEC-2A-Metalografia', -, -]"
$ws.Range("F16").Value = "[-, 'MEC-2A-Metalografia', -, -]"
